# Dash App integriert und Benni ist Dumm
# Appends the daily Germany Cases rows for 2020-05-18 .. 2020-06-03
# (rows 72-88) to the bottom of the existing table on Tabelle1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one day's row:
# Date(serial), B..Q (the 16 numeric case/death columns), in column order.
$newRows = @(
    @(43969,1652,2287,180,149,38,232,439,20,545,1512,217,149,196,54,126,139),
    @(43970,1653,2314,182,150,38,232,442,20,548,1525,222,154,198,54,128,147),
    @(43971,1672,2339,186,151,38,236,446,20,557,1537,222,154,198,54,131,149),
    @(43972,1680,2348,186,151,39,239,451,20,566,1546,224,157,201,54,134,151),
    @(43973,1680,2358,190,152,39,240,457,20,568,1547,224,157,203,54,134,151),
    @(43974,1692,2367,190,154,41,240,458,20,571,1555,226,157,204,54,136,151),
    @(43975,1697,2377,191,154,41,241,460,20,572,1563,227,157,204,54,137,152),
    @(43976,1698,2382,191,154,41,241,460,20,572,1565,228,157,205,54,137,152),
    @(43977,1707,2401,191,154,42,241,462,20,574,1571,228,158,207,54,139,153),
    @(43978,1724,2404,193,154,42,242,466,20,578,1577,229,160,208,55,140,157),
    @(43979,1724,2438,195,155,41,245,471,20,583,1582,230,161,210,55,143,158),
    @(43980,1743,2441,196,155,41,246,476,20,585,1588,230,162,209,55,144,159),
    @(43981,1748,2450,197,155,42,246,476,20,593,1595,230,163,213,55,145,161),
    @(43982,1748,2450,198,155,42,252,477,20,593,1597,230,163,213,55,145,162),
    @(43983,1752,2451,198,155,42,253,480,20,593,1598,230,163,214,55,145,162),
    @(43984,1758,2456,198,155,42,253,480,20,594,1599,230,163,212,55,145,163),
    @(43985,1760,2466,198,156,43,253,480,20,596,1607,230,163,212,55,146,166)
)

$startRow = 72
$r = $startRow
foreach ($rowData in $newRows) {
    $c = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Column A holds dates; copy the existing date number format (d-mmm-yy,
# builtin numFmtId 15) from the row above instead of leaving General,
# so the new date cells keep the same style as the rest of the column.
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $startRow + ":A" + $lastRow).NumberFormat = $ws.Range("A" + ($startRow - 1)).NumberFormat

# Restore the view: scroll so the newly added rows are visible and select
# the first empty cell right after the table (B89), matching where the
# user left off after pasting in the new data.
$ws.Range("B" + ($lastRow + 1)).Select()
